$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C12").Value = 7815
$ws.Range("C13:C23").Value = 7318
$ws.Range("C24:C148").Value = 7293
